$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8,9,10,11,13,14,15,16,21,22,23,24,30,31,32,33,37,38,39,40)
foreach ($r in $rows) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
